$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing "|-|" separators to "->" arrows, and fix content per diff
$ws.Range("C2").Value = "DM206 -> DM209"
$ws.Range("C3").Value = "DM209 -> DM223"
$ws.Range("C4").Value = "DM203 -> DM205"
$ws.Range("C6").Value = "DM205 -> DM207"
$ws.Range("C8").Value = "Glycolaldehyde -> Glycolate"
$ws.Range("C9").Value = "Glyoxylate -> tGcl Tartronate semialdehyde"
$ws.Range("C10").Value = "Hydroxypyruvate -> Glycerate"
$ws.Range("C11").Value = "tGcl Tartronate semialdehyde -> Hydroxypyruvate"
$ws.Range("D11").Value = "['Hyi']"
$ws.Range("C12").Value = "1,4-Butanediol -> 4-Hydroxy-3-Keto-Butyryl-CoA"
$ws.Range("C13").Value = "Glycolate -> Glyoxylate"
$ws.Range("C14").Value = "Glycolyl-CoA -> Glycolaldehyde"
$ws.Range("C15").Value = "Glyoxylate -> TCA cycle"
$ws.Range("C16").Value = "Acetyl-CoA -> TCA cycle"

# Add new row 17
$ws.Range("A16").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A17").Value = 10
$ws.Range("B17").Value = "PMC10269461_spectrum.04988-22-f004.jpg"
$ws.Range("C17").Value = "GCLpathway -> Acetyl-CoA"
$ws.Range("D17").Value = "[]"
$ws.Range("E17").Value = "[]"
$ws.Range("F17").Value = "[]"
$ws.Range("G17").Value = "[]"
